$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at the top; everything currently in rows 1..101 shifts
# down to rows 6..106 (the 5 trailing blank rows become 102..106).
$ws.Rows("1:5").Insert()

# Copy number formats / alignment down from the row that used to be row 1
# (now row 6) onto the freshly inserted rows so they match the rest of the
# sheet (date format on A, text format on D, F, G).
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A1:A5").PasteSpecial(-4122) | Out-Null
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D1:D5").PasteSpecial(-4122) | Out-Null
$ws.Range("F6:G6").Copy() | Out-Null
$ws.Range("F1:G5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Trailing padding on the monto strings uses U+00A0 (NBSP) x2, matching
# every other monetary shared string already in the workbook. (Use string
# interpolation, not '+', so PowerShell doesn't coerce the numeric-looking
# text into an arithmetic add.)
$nbsp = [char]0x00A0

# New transaction data (newest on top), pushing prior rows down.
$ws.Range("A1").Value2 = 41788
$ws.Range("B1").Value2 = "RETIRO ATM BP D/REINA VICTORIA"
$ws.Range("C1").Value2 = "D"
$ws.Range("D1").Value2 = "0000616755"
$ws.Range("E1").Value2 = "CENTRO DE ACOPIO NORTE"
$ws.Range("F1").Value2 = "20.00$nbsp$nbsp"
$ws.Range("G1").Value2 = "62.39"

$ws.Range("A2").Value2 = 41787
$ws.Range("B2").Value2 = "CONSUMO VISA NA SANDRY"
$ws.Range("C2").Value2 = "D"
$ws.Range("D2").Value2 = "0005149037"
$ws.Range("E2").Value2 = "INSTITUCIONAL SS.CC."
$ws.Range("F2").Value2 = "8.45$nbsp$nbsp"
$ws.Range("G2").Value2 = "82.39"

$ws.Range("A3").Value2 = 41786
$ws.Range("B3").Value2 = "RETIRO ATM BP N/GIRON 1"
$ws.Range("C3").Value2 = "D"
$ws.Range("D3").Value2 = "0000956117"
$ws.Range("E3").Value2 = "EL GIRON"
$ws.Range("F3").Value2 = "20.00$nbsp$nbsp"
$ws.Range("G3").Value2 = "90.84"

$ws.Range("A4").Value2 = 41785
$ws.Range("B4").Value2 = "CONSUMO DATA AKI MOLINEROS 161"
$ws.Range("C4").Value2 = "D"
$ws.Range("D4").Value2 = "0013346747"
$ws.Range("E4").Value2 = "INSTITUCIONAL SS.CC."
$ws.Range("F4").Value2 = "7.57$nbsp$nbsp"
$ws.Range("G4").Value2 = "110.84"

$ws.Range("A5").Value2 = 41785
$ws.Range("B5").Value2 = "CONSUMO VISA NA MENESTRAS DEL NEGRO M0"
$ws.Range("C5").Value2 = "D"
$ws.Range("D5").Value2 = "0009362885"
$ws.Range("E5").Value2 = "INSTITUCIONAL SS.CC."
$ws.Range("F5").Value2 = "11.42$nbsp$nbsp"
$ws.Range("G5").Value2 = "118.41"

# Re-create the shared "export row" formula over the new top block H1:H5
# (it used to cover H1:H14; the row that held the old master formula is
# now row 6 and no longer carries it).
$formula = "=CONCATENATE(""array('mo_fecha' => new \DateTime('"",TEXT(A1,""yyyy-mm-dd""),""'), 'mo_concepto' => '"",B1,""', 'mo_tipo' => '"",C1,""', 'mo_documento' => '"",D1,""', 'mo_oficina' => '"",E1,""', 'mo_monto' => "",F1,"", 'mo_saldo' => "",G1,"", 'mo_fecha_crea' => new \DateTime('"",TEXT(NOW(),""yyyy-mm-dd H:m:s""),""'), 'mo_quien_crea' => 1, 'mo_fecha_modifica' => NULL, 'mo_quien_modifica' => NULL, 'mo_fecha_borrado' => NULL, 'mo_quien_borra' => NULL, 'mo_borrado_logico' => false),"")"
$ws.Range("H1:H5").Formula = $formula
$ws.Range("H6").ClearContents()

# Match the saved selection in the source file.
$ws.Activate()
$ws.Range("H1:H5").Select()
